# "added task 14 and task 15 folders"
#
# DAILY-TRACKER sheet: fill in the previously-blank Task 14 detail column,
# and add the Task 15 row underneath it.
#
# TASKS-SHEET sheet: fix a typo (stray trailing backtick) in the
# "View Selector & Column Selector" task-14 entry.

$wb = $excel.ActiveWorkbook

$tracker = $wb.Worksheets.Item("DAILY-TRACKER")
$tasks   = $wb.Worksheets.Item("TASKS-SHEET")

# ---- DAILY-TRACKER: row 26 (TASK 14) ------------------------------------
# Column C already holds "Report properties" but actually documents the
# "Column selector & View selector" task; Column D was empty and now
# records that TASK 14 was about report properties.
$tracker.Range("C26").Value = "Column selector & View selector"
$tracker.Range("D26").Value = "Report properties"
$tracker.Range("B26").NumberFormat = "dd/mmm"

# ---- DAILY-TRACKER: row 27 (new, TASK 15) -------------------------------
$tracker.Range("A27").Value = 15
$tracker.Range("B27").Value = 44390
$tracker.Range("B27").NumberFormat = "dd/mmm"
$tracker.Range("C27").Value = "Report properties "
$tracker.Range("D27").Value = "Report properties Final Documentation"
$tracker.Range("E27").Value = "yes"

# ---- TASKS-SHEET: fix stray backtick in the Task 14 description --------
$tasks.Range("B71").Value = "View Selector & Column Selector"
